# Applies the "samples from 24 to 29 degrees added" edit to BermudaSalinity.
# Fills in Salinity (C), Date Run (D) and Run By (E) for a batch of rows that
# previously had blank measurements, fixes a bad "3/39/2018" text date that
# had been typed into column D for several rows (replaced with the real
# numeric date 3/29/2018 = serial 43188), and updates a couple of rows whose
# values were revised.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row groups
# ---------------------------------------------------------------------

# Rows that had completely blank C/D/E cells and now get a Salinity reading,
# a Date Run of 4/2/2018 (serial 43192) and Run By = "JM".
$newRows = @{
    36  = 37.950000000000003
    37  = 37.950000000000003
    38  = 37.94
    40  = 37.96
    41  = 37.909999999999997
    43  = 37.92
    44  = 37.96
    46  = 37.909999999999997
    47  = 37.94
    49  = 37.96
    50  = 37.93
    51  = 37.96
    52  = 37.93
    208 = 37.729999999999997
    212 = 37.74
    214 = 37.74
    215 = 37.700000000000003
    219 = 37.76
    223 = 37.74
    224 = 37.770000000000003
    225 = 37.799999999999997
    288 = 38.119999999999997
    293 = 38.1
    296 = 38.090000000000003
    299 = 37.83
    301 = 38.020000000000003
    302 = 37.979999999999997
}

# Rows whose Date Run was mistakenly entered as the text "3/39/2018" (not a
# valid date) - correct it to the real numeric date, 3/29/2018 (serial 43188).
$badDateRows = @(206, 207, 209, 210, 213, 217, 218, 221, 222, 226, 227, 228, 229, 230, 231)

# Row 220 is like $badDateRows but its Date Run cell was already bold
# (it shares the bold-header style), so the fixed cell needs a bold+date
# style instead of the plain date style used everywhere else.
$boldBadDateRow = 220

# A reference cell that already carries the plain "date" style (numFmtId 14
# + border, no bold) so we can clone that style via copy/paste instead of
# re-deriving a number format from scratch (which would create a brand new
# custom numFmt/style entry instead of reusing the existing one).
$dateStyleSource = $ws.Cells.Item(172, 4)

# ---------------------------------------------------------------------
# Fill the previously-empty rows
# ---------------------------------------------------------------------

$dateStyleSource.Copy() | Out-Null
foreach ($row in $newRows.Keys) {
    $ws.Cells.Item($row, 4).PasteSpecial(-4122) | Out-Null
}
foreach ($row in $badDateRows) {
    $ws.Cells.Item($row, 4).PasteSpecial(-4122) | Out-Null
}
$ws.Cells.Item(211, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

foreach ($row in $newRows.Keys) {
    $ws.Cells.Item($row, 3).Value2 = $newRows[$row]
    $ws.Cells.Item($row, 4).Value2 = 43192
    $ws.Cells.Item($row, 5).Value2 = "JM"
}

# ---------------------------------------------------------------------
# Fix the bad "3/39/2018" text dates -> real date serial 43188
# ---------------------------------------------------------------------

foreach ($row in $badDateRows) {
    $ws.Cells.Item($row, 4).Value2 = 43188
}

# Row 220's Date Run cell needs the bold variant of the date style. Make it
# bold first, then apply the (non-bold) date number format - doing it in
# this order reuses the existing built-in date numFmtId (14) instead of
# minting a new custom number format.
$ws.Cells.Item($boldBadDateRow, 4).Font.Bold = $true
$ws.Cells.Item($boldBadDateRow, 4).NumberFormat = "mm-dd-yy"
$ws.Cells.Item($boldBadDateRow, 4).Value2 = 43188

# ---------------------------------------------------------------------
# Row-specific corrections
# ---------------------------------------------------------------------

# Row 211: salinity reading revised, and its bad text date (style already
# fixed up above, alongside $badDateRows) corrected to 4/2/2018 (serial
# 43192) rather than 43188.
$ws.Cells.Item(211, 3).Value2 = 37.72
$ws.Cells.Item(211, 4).Value2 = 43192

# Row 309: salinity reading revised and Date Run moved from 3/28/2018
# (43187) to 4/2/2018 (43192). It already has the correct date style.
$ws.Cells.Item(309, 3).Value2 = 37.950000000000003
$ws.Cells.Item(309, 4).Value2 = 43192

# ---------------------------------------------------------------------
# View state (best effort - scroll position / selection)
# ---------------------------------------------------------------------

$ws.Activate()
$ws.Range("G212").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 204
$excel.ActiveWindow.ScrollColumn = 1
